$d = $word.ActiveDocument

function Replace-Literal($doc, $oldText, $newText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
        return $null
    }
    $rng.Text = $newText
    return $rng
}

function Delete-ParagraphContaining($doc, $text) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        Write-Output "NOT FOUND PARA: $text"
        return
    }
    $para = $rng.Paragraphs(1).Range
    $para.Delete()
}

# 1. Intro paragraph: remove Visual Studio 2012 / Windows SDK 8.0 / Windows Vista references
Replace-Literal $d "This is the DirectX SDK's Direct3D 11 sample updated to use Visual Studio 2012 and the Windows SDK 8.0 without any dependencies on legacy DirectX SDK content. This sample is a Win32 desktop DirectX 11.0 application for Windows 10, Windows 8.1, Windows 8, Windows 7, and Windows Vista Service Pack 2 with the DirectX 11.0 runtime. " "This is the DirectX SDK's Direct3D 11 sample updated to use the Windows 10 SDK without any dependencies on legacy DirectX SDK content. This sample is a Win32 desktop DirectX 11.0 application for Windows 10, Windows 8.1, Windows 8, and Windows 7. " | Out-Null

# 2. Dependencies section: Windows 8.x SDK / Windows Vista -> Windows 10 SDK / Windows 7, drop D3DCompile_46
Replace-Literal $d "When using the Windows 8.x SDK and targeting Windows Vista or later, you can include the D3DCompile_46 or D3DCompile_47 DLL side-by-side with your application copying the file from the REDIST folder. " "When using the Windows 10 SDK and targeting Windows 7 or later, you can include the D3DCompile_47 DLL side-by-side with your application copying the file from the REDIST folder. " | Out-Null

# 3. Windows kits path: 8.0 -> 10, drop "arm, "
Replace-Literal $d "(x86)%\Windows kits\8.0\" "(x86)%\Windows kits\10\" | Out-Null
Replace-Literal $d "\D3D\arm, x86 or x64" "\D3D\ x86 or x64" | Out-Null

# 4. Remove the Windows kits 8.1 and duplicate Windows kits 10 paragraphs entirely
Delete-ParagraphContaining $d "Windows kits\8.1" | Out-Null
Delete-ParagraphContaining $d "Windows kits\10" | Out-Null
